# Regenerate orders with updated distance/sizes.
# The experiment's distance codes and one of the size codes were
# renumbered:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These codes appear embedded inside many strings across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size columns),
# so do a global find/replace across the whole used range. Replace the
# longer/more-specific distance codes before anything that could be a
# substring of another token, and do the size change last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
